$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- I. UBICACIÓN GEOGRÁFICA (rows 6-11, column C values) ---
$ws.Range("C6").Value = "ewrsdfserh"
$ws.Range("C7").Value = "vgjbjkb"
$ws.Range("C8").Value = "jkjbkjkb"
$ws.Range("C9").Value = "lkklh"
$ws.Range("C10").Value = "klj"
$ws.Range("C11").Value = "kjljkl"

# --- II. INFORMACIÓN GENERAL (rows 13-20) ---
$ws.Range("C13").Value = "ñljjñllj"
$ws.Range("C14").Value = "ÑLJÑJLLÑJ"
$ws.Range("C15").Value = "JLÑJLÑ"
$ws.Range("C16").Value = "JLÑLJÑ"
$ws.Range("H16").Value = 14
$ws.Range("C17").Value = "JLÑLJÑ"
$ws.Range("H17").Value = 14
$ws.Range("C18").Value = "JLÑLJ"
$ws.Range("H18").Value = 14
$ws.Range("C19").Value = "ÑJLÑ"
$ws.Range("H19").Value = 14
$ws.Range("C20").Value = "LJÑLJÑ"

# --- IV. REGISTRO TÉCNICO DE INSTALACIÓN (rows 23-25) ---
$ws.Range("C23").Value = "LJLÑJLJÑ"
$ws.Range("G23").Value = "JLLJÑÑLJ"
$ws.Range("J23").Value = "JLÑJLÑ"
$ws.Range("C24").Value = "JLJLÑ"
$ws.Range("G24").Value = "LJÑ"
$ws.Range("J24").Value = "LJÑ"
$ws.Range("C25").Value = "JLÑLÑJ"
$ws.Range("G25").Value = "KJ"
$ws.Range("J25").Value = "LJÑ"

# --- V./VI. REGISTRO TECNICO DE FUNCIONAMIENTO / CLASIFICACION (rows 27-31) ---
$ws.Range("H27").Value = "ÑLJJK"
$ws.Range("H28").Value = "IHIH"
$ws.Range("H29").Value = "LKHLKHKL"
# C30 and H30/H31 keep their original text look; C30 holds a numeric-looking
# string so force it to stay text with a leading apostrophe (quote prefix)
$ws.Range("C30").Value = "'12"
$ws.Range("H30").Value = "KHLHKHKL"
$ws.Range("H31").Value = "KHLKHL"

# --- IIX. PERIODICIDAD DEL MANTENIMIENTO (row 34) ---
$ws.Range("E34").Value = "TRIMESTRAL"
# H34 also holds a numeric-looking string; keep it text
$ws.Range("H34").Value = "'2"

# --- X. ACCESORIOS (rows 36-38) ---
$ws.Range("A36").Value = "KÑJHLK"
# E36/E37/E38 hold numeric-looking strings ("4. ", "5. ", "6."); keep them text
$ws.Range("E36").Value = "'4. "
$ws.Range("A37").Value = "ÑLKÑJLÑJ"
$ws.Range("E37").Value = "'5. "
$ws.Range("A38").Value = "LJKLÑJLÑÑJL"
$ws.Range("E38").Value = "'6."
